# Update workbook with Daria's latest UK parametric-matching estimates.
# - Sheet1 "UK" becomes "Info": a short description/author/edit-date block.
# - Sheet2 "IT" becomes "Parameters": refreshed UK regression estimates.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "UK" -> "Info" ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Info"

# Drop the old parameter rows 4-6 (var_dag_diff / var_wage_diff / cov_dag_wage_diff)
# and the numeric estimate in rows 2-3, col B, leaving a 3-row info block.
$ws1.Rows("4:6").Delete()
$ws1.Range("B2:B3").ClearContents()

$ws1.Range("A1").Value = "Description:"
$ws1.Range("B1").Value = "Estimates for the parametric couple matching process"
$ws1.Range("A2").Value = "Authors: Patryk Bronka, Daria Popova"
$ws1.Range("A3").Value = "Last edit: 4 July 2025 DP"

$ws1.Columns("A").ColumnWidth = 38.86
[void]$ws1.Range("A10").Select()

# --- Sheet 2: "IT" -> "Parameters" ------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Parameters"

$ws2.Range("A1").Value = "Parameter"
$ws2.Range("B1").Value = "Value"

$ws2.Range("A2").Value = "mean_dag_diff"
$ws2.Range("B2").Value = 1.9909832261486442

$ws2.Range("A3").Value = "mean_wage_diff"
$ws2.Range("B3").Value = 3.5943392901166864

$ws2.Range("A4").Value = "var_dag_diff"
$ws2.Range("B4").Value = 22.366846155599802

$ws2.Range("A5").Value = "var_wage_diff"
$ws2.Range("B5").Value = 126.41297991430548

$ws2.Range("A6").Value = "cov_dag_wage_diff"
$ws2.Range("B6").Value = 0.99009364374758324
